$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers need to be
# forced to text format first, so Excel keeps them as literal strings
# (preserving exact formatting such as trailing zeros) instead of parsing
# them into numeric values.
$textCells = @('D6','D7','D8','D9','D10','D11','D14','D15','D17','D18','D19','D21','D22','D23','D26','D27','D28','D29','D30','D31','D33','D34','D37','D39','D40','D41','D42','D44','D45','D46','D47','D48','D49','D50','D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '22.373.21'
$ws.Range('E2').Value = '  -4.51%  '
$ws.Range('D3').Value = '1.567.06'
$ws.Range('E3').Value = '  -4.75%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').Value = '290.27'
$ws.Range('E6').Value = '  -2.94%  '
$ws.Range('D7').Value = '0.3671'
$ws.Range('E7').Value = '  -3.06%  '
$ws.Range('D8').Value = '49.53'
$ws.Range('E8').Value = '  -0.48%  '
$ws.Range('D9').Value = '0.3387'
$ws.Range('E9').Value = '  -3.42%  '
$ws.Range('D10').Value = '1.170'
$ws.Range('E10').Value = '  -3.47%  '
$ws.Range('D11').Value = '0.07609'
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('E13').Value = '  -3.98%  '
$ws.Range('D14').Value = '6.052'
$ws.Range('E14').Value = '  -4.87%  '
$ws.Range('D15').Value = '6.893'
$ws.Range('E15').Value = '  -5.77%  '
$ws.Range('D16').Value = '1.576.90'
$ws.Range('E16').Value = '  -3.91%  '
$ws.Range('D17').Value = '0.00001134'
$ws.Range('E17').Value = '  -5.63%  '
$ws.Range('D18').Value = '89.08'
$ws.Range('E18').Value = '  -7.72%  '
$ws.Range('D19').Value = '0.06762'
$ws.Range('E19').Value = '  -3.58%  '
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('D21').Value = '6.227'
$ws.Range('E21').Value = '  -7.30%  '
$ws.Range('D22').Value = '0.5341'
$ws.Range('E22').Value = '  -6.82%  '
$ws.Range('D23').Value = '16.50'
$ws.Range('E23').Value = '  -4.85%  '
$ws.Range('E24').Value = '  -2.59%  '
$ws.Range('D25').Value = '22.380.58'
$ws.Range('E25').Value = '  -4.52%  '
$ws.Range('D26').Value = '2.385'
$ws.Range('E26').Value = '  -4.56%  '
$ws.Range('D27').Value = '2.980'
$ws.Range('E27').Value = '  +2.45%  '
$ws.Range('D28').Value = '19.90'
$ws.Range('E28').Value = '  -4.53%  '
$ws.Range('D29').Value = '145.71'
$ws.Range('E29').Value = '  -4.81%  '
$ws.Range('D30').Value = '4.960'
$ws.Range('E30').Value = '  -4.70%  '
$ws.Range('D31').Value = '125.40'
$ws.Range('E31').Value = '  -5.38%  '
$ws.Range('D32').Value = '1.755.53'
$ws.Range('E32').Value = '  -3.96%  '
$ws.Range('D33').Value = '1.036'
$ws.Range('E33').Value = '  +5.77%  '
$ws.Range('D34').Value = '6.249'
$ws.Range('E34').Value = '  -9.07%  '
$ws.Range('E35').Value = '  -6.05%  '
$ws.Range('E36').Value = '  -9.58%  '
$ws.Range('D37').Value = '0.08460'
$ws.Range('E37').Value = '  -3.20%  '
$ws.Range('E38').Value = '  -5.87%  '
$ws.Range('D39').Value = '0.2327'
$ws.Range('E39').Value = '  -4.14%  '
$ws.Range('D40').Value = '0.06554'
$ws.Range('E40').Value = '  -3.81%  '
$ws.Range('D41').Value = '5.529'
$ws.Range('E41').Value = '  -6.25%  '
$ws.Range('D42').Value = '11.82'
$ws.Range('E42').Value = '  -8.00%  '
$ws.Range('E43').Value = '  -3.94%  '
$ws.Range('D44').Value = '0.6371'
$ws.Range('E44').Value = '  -7.16%  '
$ws.Range('D45').Value = '14.38'
$ws.Range('E45').Value = '  -8.21%  '
$ws.Range('D46').Value = '1.000'
$ws.Range('D47').Value = '0.6005'
$ws.Range('E47').Value = '  -5.19%  '
$ws.Range('D48').Value = '3.780'
$ws.Range('E48').Value = '  -3.20%  '
$ws.Range('D49').Value = '2.131'
$ws.Range('E49').Value = '  -5.28%  '
$ws.Range('D50').Value = '1.263'
$ws.Range('E50').Value = '  +7.36%  '
$ws.Range('D51').Value = '123.24'
$ws.Range('E51').Value = '  -2.85%  '
